$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.138256430625916
$ws.Range("B1").Value = 2.608568429946899
$ws.Range("C1").Value = 6.981382369995117
$ws.Range("D1").Value = 2.075896501541138
$ws.Range("E1").Value = 1.189907789230347
